# Fix problems with yellow and red rows:
# - row 6 (port 3) gets its missing B/C/D details filled in
# - row 4's ALP_6520_DISK_1..72 entries get a trailing "." appended
# - selection moves to D4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-missing details for row 6 first, so that the
# shared-string table picks up these new strings before the row 4 edit
# (matches the order new unique strings were introduced upstream).
$ws.Range("B6").Value = "ALP_X6_DISK_1.F_PORT.5."
$ws.Range("D6").Value = "ALP_X6_DISK_2.F_PORT.5."
$ws.Range("C6").Value = "port 3"

# Correct the row 4 disk identifiers by appending a trailing period.
$ws.Range("B4").Value = "ALP_6520_DISK_1..72."
$ws.Range("D4").Value = "ALP_6520_DISK_1..72."

# Update the active selection to D4.
$ws.Range("D4").Select()
